# ============================================================
# Add files via upload
# Added some back-end validations
# ============================================================

$wb = $excel.ActiveWorkbook

$wsFront = $wb.Worksheets.Item(1)
$wsBack  = $wb.Worksheets.Item(2)

# ---- Rename sheets --------------------------------------------------
$wsFront.Name = "Front-End"
$wsBack.Name  = "Back- End"

# ---- Front-End sheet: remove the stray spacer row 49 ----------------
$wsFront.Range("B49:E49").Clear()

# ---- Back- End sheet: remove the old placeholder content -------------
$wsBack.Range("E3:H3").UnMerge()
$wsBack.Range("E4:H4").UnMerge()
$wsBack.Range("E3:H4").Clear()

# ---- Back- End sheet: column widths ----------------------------------
$wsBack.Columns.Item(6).ColumnWidth = 41.77734375
$wsBack.Columns.Item(7).ColumnWidth = 69.88671875

# ---- Section: Customer Entity ----------------------------------------
$wsBack.Range("E3").Value = "Section: Customer Entity"
$r = $wsBack.Range("E3:G3")
$r.Merge()
$r.Interior.ColorIndex = 41
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108

$wsBack.Range("E4").Value = "SR No"
$wsBack.Range("F4").Value = "Test Cases"
$wsBack.Range("G4").Value = "Result"
$r = $wsBack.Range("E4:G4")
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108

$wsBack.Range("E5").Value = 1
$wsBack.Range("F5").Value = "Contains other than customer email format"
$wsBack.Range("G5").Value = "Customer email contains characters, @ and digits with proper format"
$wsBack.Range("E5").HorizontalAlignment = -4108
$wsBack.Range("E5").VerticalAlignment = -4108
$wsBack.Range("F5:G5").VerticalAlignment = -4108

$wsBack.Range("E6").Value = 2
$wsBack.Range("F6").Value = "Contains other than customer phone no. format"
$wsBack.Range("G6").Value = "Customer Phone no. must contain digits and its only 10 digits"
$wsBack.Range("E6").HorizontalAlignment = -4108
$wsBack.Range("F6:G6").VerticalAlignment = -4108

$wsBack.Range("E7").Value = 3
$wsBack.Range("F7").Value = "Contains other than customer adhaar no. format"
$wsBack.Range("G7").Value = "Customer Adhaar no. must contain digits and its only 12 digits with proper format"
$wsBack.Range("E7").HorizontalAlignment = -4108
$wsBack.Range("F7:G7").VerticalAlignment = -4108

# ---- Section: Card Entity ---------------------------------------------
$wsBack.Range("E9").Value = "Section: Card Entity"
$r = $wsBack.Range("E9:G9")
$r.Merge()
$r.Interior.ColorIndex = 41
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108

$wsBack.Range("E10").Value = "SR No"
$wsBack.Range("F10").Value = "Test Cases"
$wsBack.Range("G10").Value = "Result"
$r = $wsBack.Range("E10:G10")
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4108

$wsBack.Range("E11").Value = 1
$wsBack.Range("F11").Value = "Card Cvv not contain three digits"
$wsBack.Range("G11").Value = "Card Cvv must cotain three digits only"
$wsBack.Range("E11").HorizontalAlignment = -4108
$wsBack.Range("E11").VerticalAlignment = -4108

# ---- View state --------------------------------------------------------
$wsFront.Activate()
$excel.ActiveWindow.Zoom = 109
$wsFront.Range("D14").Select()

$wsBack.Activate()
$excel.ActiveWindow.Zoom = 110
$wsBack.Range("I5").Select()

Write-Host "done"
